$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.746.59'
$ws.Range('E2').Value = '  -3.48%  '

$ws.Range('D3').Value = '2.908.00'
$ws.Range('E3').Value = '  -4.20%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = '582.55'
$ws.Range('E5').Value = '  -1.91%  '

$ws.Range('D6').Value = '144.04'
$ws.Range('E6').Value = '  -6.61%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('E8').Value = '  -2.91%  '

$ws.Range('D9').Value = '2.906.74'
$ws.Range('E9').Value = '  -4.07%  '

$ws.Range('D10').Value = '6.79'
$ws.Range('E10').Value = '  +5.39%  '

$ws.Range('E11').Value = '  -4.85%  '

$ws.Range('E12').Value = '  -4.54%  '

$ws.Range('E13').Value = '  -4.02%  '

$ws.Range('E14').Value = '  -6.37%  '

$ws.Range('E15').Value = '  -0.03%  '

$ws.Range('D16').Value = '3.391.08'
$ws.Range('E16').Value = '  -4.10%  '

$ws.Range('D17').Value = '60.713.14'
$ws.Range('E17').Value = '  -3.49%  '

$ws.Range('D18').Value = '6.74'
$ws.Range('E18').Value = '  -5.33%  '

$ws.Range('D19').Value = '2.910.31'
$ws.Range('E19').Value = '  -4.01%  '

$ws.Range('D20').Value = '430.89'
$ws.Range('E20').Value = '  -5.03%  '

$ws.Range('D21').Value = '13.59'
$ws.Range('E21').Value = '  -5.08%  '

$ws.Range('D22').Value = '0.682'
$ws.Range('E22').Value = '  -2.30%  '

$ws.Range('E23').Value = '  -4.94%  '

$ws.Range('D24').Value = '80.32'
$ws.Range('E24').Value = '  -3.39%  '

$ws.Range('D25').Value = '10.92'
$ws.Range('E25').Value = '  -2.46%  '

$ws.Range('E26').Value = '  -5.07%  '

$ws.Range('E27').Value = '  -4.23%  '

$ws.Range('E28').Value = '  -0.02%  '

$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.01%  '

$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = '7.16'
$ws.Range('E30').Value = '  -3.98%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.61'
$ws.Range('E31').Value = '  -3.16%  '

$ws.Range('E32').Value = '  -2.25%  '

$ws.Range('D33').Value = '26.46'
$ws.Range('E33').Value = '  -4.12%  '

$ws.Range('E34').Value = '  -4.42%  '

$ws.Range('D35').Value = '0.0₃0865'
$ws.Range('E35').Value = '  -0.71%  '

$ws.Range('E36').Value = '  -3.10%  '

$ws.Range('D37').Value = '5.64'
$ws.Range('E37').Value = '  -4.97%  '

$ws.Range('E38').Value = '  -6.21%  '

$ws.Range('D39').Value = '0.127'
$ws.Range('E39').Value = '  -0.46%  '

$ws.Range('D40').Value = '49.72'
$ws.Range('E40').Value = '  -1.82%  '

$ws.Range('E41').Value = '  -5.74%  '

$ws.Range('D42').Value = '8.65'
$ws.Range('E42').Value = '  -4.75%  '

$ws.Range('D43').Value = '0.294'
$ws.Range('E43').Value = '  -4.57%  '

$ws.Range('D44').Value = '41.39'
$ws.Range('E44').Value = '  -1.23%  '

$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Value = '374.78'
$ws.Range('E45').Value = '  -5.82%  '

$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0347'
$ws.Range('E46').Value = '  -3.41%  '

$ws.Range('D47').Value = '2.667.22'
$ws.Range('E47').Value = '  -2.37%  '

$ws.Range('D48').Value = '132.71'
$ws.Range('E48').Value = '  +0.25%  '

$ws.Range('E49').Value = '  -0.04%  '

$ws.Range('D50').Value = '24.30'
$ws.Range('E50').Value = '  -0.79%  '

$ws.Range('D51').Value = '0.106'
$ws.Range('E51').Value = '  -2.16%  '
